$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-07-12 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-13 Sunday", 2) | Out-Null

# Update table cells by position (handles duplicate source text correctly)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "25+52=77"
$tbl.Cell(1, 2).Range.Text = "26+54=80"
$tbl.Cell(1, 3).Range.Text = "99-80=19"
$tbl.Cell(1, 4).Range.Text = "0+35=35"
$tbl.Cell(1, 5).Range.Text = "34-25=9"
$tbl.Cell(2, 1).Range.Text = "5+16=21"
$tbl.Cell(2, 2).Range.Text = "24+1=25"
$tbl.Cell(2, 3).Range.Text = "24-23=1"
$tbl.Cell(2, 4).Range.Text = "1+19=20"
$tbl.Cell(2, 5).Range.Text = "0+55=55"
$tbl.Cell(3, 1).Range.Text = "55-38=17"
$tbl.Cell(3, 2).Range.Text = "45+14=59"
$tbl.Cell(3, 3).Range.Text = "66-15=51"
$tbl.Cell(3, 4).Range.Text = "93-29=64"
$tbl.Cell(3, 5).Range.Text = "66+15=81"
$tbl.Cell(4, 1).Range.Text = "0+59=59"
$tbl.Cell(4, 2).Range.Text = "26+8=34"
$tbl.Cell(4, 3).Range.Text = "4+90=94"
$tbl.Cell(4, 4).Range.Text = "22+38=60"
$tbl.Cell(4, 5).Range.Text = "95+2=97"
$tbl.Cell(5, 1).Range.Text = "16+39=55"
$tbl.Cell(5, 2).Range.Text = "43-39=4"
$tbl.Cell(5, 3).Range.Text = "20+53=73"
$tbl.Cell(5, 4).Range.Text = "4+4=8"
$tbl.Cell(5, 5).Range.Text = "92-37=55"
$tbl.Cell(6, 1).Range.Text = "63-45=18"
$tbl.Cell(6, 2).Range.Text = "29+28=57"
$tbl.Cell(6, 3).Range.Text = "30-1=29"
$tbl.Cell(6, 4).Range.Text = "50-28=22"
$tbl.Cell(6, 5).Range.Text = "18+73=91"
$tbl.Cell(7, 1).Range.Text = "89-7=82"
$tbl.Cell(7, 2).Range.Text = "60+15=75"
$tbl.Cell(7, 3).Range.Text = "3+73=76"
$tbl.Cell(7, 4).Range.Text = "33+58=91"
$tbl.Cell(7, 5).Range.Text = "46-23=23"
$tbl.Cell(8, 1).Range.Text = "3+80=83"
$tbl.Cell(8, 2).Range.Text = "42+35=77"
$tbl.Cell(8, 3).Range.Text = "10+34=44"
$tbl.Cell(8, 4).Range.Text = "80-35=45"
$tbl.Cell(8, 5).Range.Text = "53+17=70"
$tbl.Cell(9, 1).Range.Text = "66-24=42"
$tbl.Cell(9, 2).Range.Text = "4+86=90"
$tbl.Cell(9, 3).Range.Text = "90+5=95"
$tbl.Cell(9, 4).Range.Text = "28+59=87"
$tbl.Cell(9, 5).Range.Text = "91-48=43"
$tbl.Cell(10, 1).Range.Text = "7+79=86"
$tbl.Cell(10, 2).Range.Text = "57+42=99"
$tbl.Cell(10, 3).Range.Text = "58-39=19"
$tbl.Cell(10, 4).Range.Text = "61+1=62"
$tbl.Cell(10, 5).Range.Text = "3+25=28"
$tbl.Cell(11, 1).Range.Text = "8+17=25"
$tbl.Cell(11, 2).Range.Text = "15+13=28"
$tbl.Cell(11, 3).Range.Text = "52+18=70"
$tbl.Cell(11, 4).Range.Text = "94+1=95"
$tbl.Cell(11, 5).Range.Text = "27+42=69"
$tbl.Cell(12, 1).Range.Text = "41+32=73"
$tbl.Cell(12, 2).Range.Text = "69-59=10"
$tbl.Cell(12, 3).Range.Text = "64-22=42"
$tbl.Cell(12, 4).Range.Text = "1+75=76"
$tbl.Cell(12, 5).Range.Text = "72-53=19"
$tbl.Cell(13, 1).Range.Text = "27+55=82"
$tbl.Cell(13, 2).Range.Text = "79-0=79"
$tbl.Cell(13, 3).Range.Text = "71+3=74"
$tbl.Cell(13, 4).Range.Text = "23+4=27"
$tbl.Cell(13, 5).Range.Text = "60-44=16"
$tbl.Cell(14, 1).Range.Text = "3+4=7"
$tbl.Cell(14, 2).Range.Text = "46+8=54"
$tbl.Cell(14, 3).Range.Text = "69-61=8"
$tbl.Cell(14, 4).Range.Text = "16+48=64"
$tbl.Cell(14, 5).Range.Text = "76-14=62"
$tbl.Cell(15, 1).Range.Text = "87-35=52"
$tbl.Cell(15, 2).Range.Text = "93+1=94"
$tbl.Cell(15, 3).Range.Text = "60+39=99"
$tbl.Cell(15, 4).Range.Text = "71-2=69"
$tbl.Cell(15, 5).Range.Text = "2+29=31"
$tbl.Cell(16, 1).Range.Text = "26+65=91"
$tbl.Cell(16, 2).Range.Text = "90-90=0"
$tbl.Cell(16, 3).Range.Text = "62+4=66"
$tbl.Cell(16, 4).Range.Text = "30-15=15"
$tbl.Cell(16, 5).Range.Text = "87-30=57"
$tbl.Cell(17, 1).Range.Text = "20+15=35"
$tbl.Cell(17, 2).Range.Text = "33+13=46"
$tbl.Cell(17, 3).Range.Text = "0+79=79"
$tbl.Cell(17, 4).Range.Text = "68-66=2"
$tbl.Cell(17, 5).Range.Text = "82-73=9"
$tbl.Cell(18, 1).Range.Text = "47+41=88"
$tbl.Cell(18, 2).Range.Text = "81-2=79"
$tbl.Cell(18, 3).Range.Text = "39+1=40"
$tbl.Cell(18, 4).Range.Text = "27+64=91"
$tbl.Cell(18, 5).Range.Text = "37+15=52"
$tbl.Cell(19, 1).Range.Text = "48+11=59"
$tbl.Cell(19, 2).Range.Text = "20+25=45"
$tbl.Cell(19, 3).Range.Text = "97-81=16"
$tbl.Cell(19, 4).Range.Text = "42-18=24"
$tbl.Cell(19, 5).Range.Text = "7+78=85"
$tbl.Cell(20, 1).Range.Text = "62-23=39"
$tbl.Cell(20, 2).Range.Text = "7+20=27"
$tbl.Cell(20, 3).Range.Text = "40+28=68"
$tbl.Cell(20, 4).Range.Text = "9+86=95"
$tbl.Cell(20, 5).Range.Text = "97-97=0"
